$d = $word.ActiveDocument

$replacements = @(
    @("293÷6=", "176÷3="),
    @("563÷2=", "192÷9="),
    @("795÷4=", "952÷8="),
    @("237÷7=", "410÷8="),
    @("697÷7=", "604÷6="),
    @("207÷4=", "896÷8="),
    @("801÷7=", "961÷5="),
    @("689÷2=", "260÷3="),
    @("713÷9=", "934÷6="),
    @("134÷7=", "299÷5="),
    @("276÷9=", "211÷5="),
    @("749÷8=", "872÷2="),
    @("852÷7=", "134÷2="),
    @("428÷6=", "705÷5="),
    @("267÷5=", "100÷2="),
    @("381÷3=", "568÷4="),
    @("112÷9=", "961÷7="),
    @("942÷5=", "797÷5="),
    @("148÷9=", "222÷3="),
    @("291÷5=", "126÷7="),
    @("688÷5=", "395÷4="),
    @("492÷7=", "525÷3="),
    @("490÷2=", "334÷4="),
    @("257÷9=", "376÷3="),
    @("941÷2=", "498÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
